$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "42÷7=6, 0" "97÷3=32, 1"
Replace-Text "75÷8=9, 3" "18÷3=6, 0"
Replace-Text "10÷2=5, 0" "45÷9=5, 0"
Replace-Text "62÷2=31, 0" "53÷6=8, 5"
Replace-Text "39÷3=13, 0" "43÷8=5, 3"

Replace-Text "54÷5=10, 4" "58÷3=19, 1"
Replace-Text "62÷6=10, 2" "14÷6=2, 2"
Replace-Text "29÷8=3, 5" "45÷4=11, 1"
Replace-Text "77÷9=8, 5" "76÷4=19, 0"
Replace-Text "25÷3=8, 1" "77÷9=8, 5"

Replace-Text "82÷2=41, 0" "38÷6=6, 2"
Replace-Text "57÷2=28, 1" "64÷3=21, 1"
Replace-Text "69÷7=9, 6" "54÷7=7, 5"
Replace-Text "39÷9=4, 3" "34÷5=6, 4"
Replace-Text "52÷9=5, 7" "34÷4=8, 2"

Replace-Text "30÷3=10, 0" "53÷3=17, 2"
Replace-Text "58÷2=29, 0" "21÷7=3, 0"
Replace-Text "91÷9=10, 1" "72÷4=18, 0"
Replace-Text "92÷6=15, 2" "57÷8=7, 1"
Replace-Text "18÷4=4, 2" "16÷9=1, 7"

Replace-Text "25÷4=6, 1" "79÷9=8, 7"
Replace-Text "33÷4=8, 1" "94÷5=18, 4"
Replace-Text "23÷9=2, 5" "34÷4=8, 2"
Replace-Text "60÷9=6, 6" "31÷3=10, 1"
Replace-Text "39÷5=7, 4" "66÷3=22, 0"
